$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8207950592041016
$ws.Range("B1").Value = 1.284064769744873
$ws.Range("C1").Value = 2.68292498588562
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.659247398376465
